# Apply corrected-error-estimation / projected-years update to SoIB_summaries.xlsx
$wb = $excel.ActiveWorkbook

# --- Sheet: "Trends Status" ---
$ws1 = $wb.Worksheets.Item("Trends Status")

$ws1.Cells.Item(2, 3).Value = 3          # C2: 2 -> 3
$ws1.Cells.Item(2, 5).Value = 9.1        # E2: 7.7 -> 9.1

$ws1.Cells.Item(3, 3).Value = 8          # C3: 6 -> 8
$ws1.Cells.Item(3, 5).Value = 24.2       # E3: 23.1 -> 24.2

$ws1.Cells.Item(4, 2).Value = 9          # B4: 8 -> 9
$ws1.Cells.Item(4, 3).Value = 18         # C4: 17 -> 18
$ws1.Cells.Item(4, 4).Value = 47.4       # D4: 44.4 -> 47.4
$ws1.Cells.Item(4, 5).Value = 54.5       # E4: 65.40000000000001 -> 54.5

$ws1.Cells.Item(5, 3).Value = 4          # C5: 1 -> 4
$ws1.Cells.Item(5, 4).Value = 31.6       # D5: 33.3 -> 31.6
$ws1.Cells.Item(5, 5).Value = 12.1       # E5: 3.8 -> 12.1

$ws1.Cells.Item(6, 4).Value = 21.1       # D6: 22.2 -> 21.1

$ws1.Cells.Item(7, 2).Value = 45         # B7: 46 -> 45
$ws1.Cells.Item(7, 3).Value = 60         # C7: 67 -> 60

# --- Sheet: "Species qualification" ---
$ws4 = $wb.Worksheets.Item("Species qualification")

$ws4.Cells.Item(3, 3).Value = 19         # C3: 18 -> 19
$ws4.Cells.Item(4, 3).Value = 33         # C4: 26 -> 33

# --- Sheet: "Interannual update - High Pri" ---
$ws5 = $wb.Worksheets.Item("Interannual update - High Pri")

$ws5.Cells.Item(2, 2).Value = 64         # B2: 63 -> 64
$ws5.Cells.Item(2, 3).Value = 62.1       # C2: 61.2 -> 62.1
$ws5.Cells.Item(2, 4).Value = 64         # D2: 63 -> 64
$ws5.Cells.Item(2, 5).Value = 76.2       # E2: 75 -> 76.2

$ws5.Cells.Item(4, 2).Value = 37         # B4: 38 -> 37
$ws5.Cells.Item(4, 3).Value = 35.9       # C4: 36.9 -> 35.9
$ws5.Cells.Item(4, 4).Value = 20         # D4: 21 -> 20
$ws5.Cells.Item(4, 5).Value = 23.8       # E4: 25 -> 23.8
